$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New data block (columns F:H, rows 2:11) driven by 01.dat
# ---------------------------------------------------------------------------

$idx = @(0,1,2,3,4,5,6,7,8,9)
for ($i = 0; $i -lt $idx.Length; $i++) {
  $r = 2 + $i
  $ws.Range("F$r").Value = $idx[$i]
}

$codes = @(
  ":001101001",
  ":011011011",
  ":001101011",
  ":011111010",
  ":011111010",
  ":111011011",
  ":010111110",
  ":111011010",
  ":011101010",
  ":111110010"
)
for ($i = 0; $i -lt $codes.Length; $i++) {
  $r = 2 + $i
  $ws.Range("G$r").Value = $codes[$i]
}

$ws.Range("H2").Formula = "=E5+E4+E7+E10"
$ws.Range("H3").Formula = "=E3+E4+E6+E7+E9+E10"
$ws.Range("H4").Formula = "=E4+E5+E7+E9+E10"
$ws.Range("H5").Formula = "=E9+E7+E6+E5+E4+E3"
$ws.Range("H6").Formula = "=E9+E7+E6+E5+E4+E3"
$ws.Range("H7").Formula = "=E2+E3+E4+E6+E7+E9+E10"
$ws.Range("H8").Formula = "=E3+E5+E6+E7+E8+E9"

# H9:H11 stay empty but still need the accounting number format applied
# (done together with H2:H11 below).

# ---------------------------------------------------------------------------
# Second little "Road / Area / PCI" table, rows 12:17
# ---------------------------------------------------------------------------

$ws.Range("A12").Value = "Road"
$ws.Range("B12").Value = "Area"
$ws.Range("C12").Value = "PCI"

$roadRows = @(
  @(1, 2178.46,  33.94861),
  @(2, 12699.86, 72.22985),
  @(3, 20431.81, 71.69824),
  @(4, 1067.78,  90.31781),
  @(5, 7431.95,  100)
)
for ($i = 0; $i -lt $roadRows.Length; $i++) {
  $r = 13 + $i
  $ws.Range("A$r").Value = $roadRows[$i][0]
  $ws.Range("B$r").Value = $roadRows[$i][1]
  $ws.Range("C$r").Value = $roadRows[$i][2]
}

# ---------------------------------------------------------------------------
# Formatting
# ---------------------------------------------------------------------------

# Header row of the little table: yellow fill + thin border + centered
$tblHeader = $ws.Range("A12:C12")
$tblHeader.VerticalAlignment = -4108
$tblHeader.HorizontalAlignment = -4108
$tblHeader.Interior.Color = 65535
$tblHeader.Borders.LineStyle = 1

# Data rows of the little table: thin border + centered, no fill
$tblBody = $ws.Range("A13:C17")
$tblBody.VerticalAlignment = -4108
$tblBody.HorizontalAlignment = -4108
$tblBody.Borders.LineStyle = 1

# Cost column formatted as currency
$ws.Range("E2:E10").NumberFormat = """$""#,##0.00"

# New "total" column formatted as accounting-style with red negatives
$ws.Range("H2:H11").HorizontalAlignment = -4108
$ws.Range("H2:H11").NumberFormat = "#,##0.00_);[Red]\(#,##0.00\)"

# Index column centered
$ws.Range("F2:F11").HorizontalAlignment = -4108

# Bit-code column centered (G2 also picks up vertical centering, matching
# the look of the adjacent road-id column)
$ws.Range("G2").HorizontalAlignment = -4108
$ws.Range("G2").VerticalAlignment = -4108
$ws.Range("G3:G11").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------------

$ws.Columns.Item(5).ColumnWidth = 14.29
$ws.Columns.Item(7).ColumnWidth = 11.99
$ws.Columns.Item(8).ColumnWidth = 14.29

# ---------------------------------------------------------------------------
# Selection / window
# ---------------------------------------------------------------------------

$ws.Range("E10").Select() | Out-Null

Write-Host "done"
